# Insert a new weekly record for "Macroferia Regional de Talca - Brócoli"
# just before the existing row 183, shifting all subsequent rows down by
# one (183->184, ..., 210->211).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 183 (pushes old 183..210 down to 184..211).
$ws.Rows(183).Insert()

# Populate the new row 183 with the new weekly entry.
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44491
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100112023
$ws.Range("G183").Value = "Brócoli"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 4000
$ws.Range("K183").Value = 600
$ws.Range("L183").Value = 600
$ws.Range("M183").Value = 600
$ws.Range("N183").Value = "$/unidad"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 600
$ws.Range("Q183").Value = 1
$ws.Range("R183").Value = "Hortaliza"
